# Add a new "Rate" format option and apply it (with provenance) to the
# infdeath / susdeath / foi parameter rows on the Parameters sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")

$ws.Range("C5").Value = "Rate"
$ws.Range("C6").Value = "Rate"
$ws.Range("C7").Value = "Rate"

$ws.Columns.Item(2).ColumnWidth = 40.333333333333336
$ws.Columns.Item(3).ColumnWidth = 18.833333333333332
$ws.Columns.Item(7).ColumnWidth = 12.5

$ws.Range("E9").Select()
